# Regenerate the lattice-multiplication exercise values in the 5x3 table.
# Each table cell holds 5 lines of text separated by manual line breaks
# (<w:br/> == Chr(11) in Word's Range.Text):
#   1) "AA x BB"      - the two factors being multiplied
#   2) "  C    D"     - the individual digits of the second factor, spaced out
#   3) "  ----"       - a fixed separator (unchanged by this edit)
#   4) "E|    |"      - first lattice row stub
#   5) "F|    |"      - second lattice row stub
#
# The cell/row/column layout of the table itself does not change - only the
# text content of each cell is regenerated with new values.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$NL = [char]11

# New content for every cell, keyed by (row, col) with 1-based indices.
$cellData = @(
    @{ Row = 1; Col = 1; Lines = @("23 x 40", "  4    0", "  ----", "2|    |", "3|    |") }
    @{ Row = 1; Col = 2; Lines = @("91 x 97", "  9    7", "  ----", "9|    |", "1|    |") }
    @{ Row = 1; Col = 3; Lines = @("61 x 23", "  2    3", "  ----", "6|    |", "1|    |") }

    @{ Row = 2; Col = 1; Lines = @("12 x 81", "  8    1", "  ----", "1|    |", "2|    |") }
    @{ Row = 2; Col = 2; Lines = @("42 x 81", "  8    1", "  ----", "4|    |", "2|    |") }
    @{ Row = 2; Col = 3; Lines = @("51 x 36", "  3    6", "  ----", "5|    |", "1|    |") }

    @{ Row = 3; Col = 1; Lines = @("98 x 90", "  9    0", "  ----", "9|    |", "8|    |") }
    @{ Row = 3; Col = 2; Lines = @("26 x 67", "  6    7", "  ----", "2|    |", "6|    |") }
    @{ Row = 3; Col = 3; Lines = @("67 x 44", "  4    4", "  ----", "6|    |", "7|    |") }

    @{ Row = 4; Col = 1; Lines = @("91 x 75", "  7    5", "  ----", "9|    |", "1|    |") }
    @{ Row = 4; Col = 2; Lines = @("95 x 67", "  6    7", "  ----", "9|    |", "5|    |") }
    @{ Row = 4; Col = 3; Lines = @("88 x 98", "  9    8", "  ----", "8|    |", "8|    |") }

    @{ Row = 5; Col = 1; Lines = @("77 x 73", "  7    3", "  ----", "7|    |", "7|    |") }
    @{ Row = 5; Col = 2; Lines = @("54 x 39", "  3    9", "  ----", "5|    |", "4|    |") }
    @{ Row = 5; Col = 3; Lines = @("62 x 14", "  1    4", "  ----", "6|    |", "2|    |") }
)

foreach ($entry in $cellData) {
    $cell = $t.Cell($entry.Row, $entry.Col)
    $newText = [string]::Join($NL, $entry.Lines)
    $cell.Range.Text = $newText
}

Write-Output "Updated $($cellData.Count) cells"
